$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-20 00:52:53"
$wsZhCn.Range("P3").Value = ""

$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-20 00:52:59"
$wsDeDe.Range("P3").Value = ""

$wsZhCn.Columns.Item(16).ColumnWidth = 12.9
$wsDeDe.Columns.Item(16).ColumnWidth = 12.9
